$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated figures for 2022-02-24 data refresh (Fonds de solidarite volet 1)
$ws.Range("C2").Value2 = 766207
$ws.Range("E2").Value2 = 1428874112

$ws.Range("C10").Value2 = 345453
$ws.Range("E10").Value2 = 1816729395

$ws.Range("C13").Value2 = 187763
$ws.Range("E13").Value2 = 1162522376

$ws.Range("C36").Value2 = 211191
$ws.Range("E36").Value2 = 404207239

$ws.Range("C57").Value2 = 31590
$ws.Range("E57").Value2 = 162173592

$ws.Range("C67").Value2 = 27092
$ws.Range("E67").Value2 = 168553916

$ws.Range("C72").Value2 = 331302
$ws.Range("E72").Value2 = 635354885

$ws.Range("C78").Value2 = 178406
$ws.Range("E78").Value2 = 892017619

$ws.Range("C79").Value2 = 680
$ws.Range("E79").Value2 = 20349120

$ws.Range("C91").Value2 = 18405
$ws.Range("E91").Value2 = 72117092

$ws.Range("C93").Value2 = 16555
$ws.Range("E93").Value2 = 48245887

$ws.Range("C112").Value2 = 145188
$ws.Range("E112").Value2 = 715584525

$ws.Range("C115").Value2 = 81784
$ws.Range("D115").Value2 = 14448
$ws.Range("E115").Value2 = 435971242

$ws.Range("C121").Value2 = 1305809
$ws.Range("E121").Value2 = 2273436091

$ws.Range("C127").Value2 = 9137
$ws.Range("E127").Value2 = 110235394

$ws.Range("C128").Value2 = 280
$ws.Range("E128").Value2 = 5719119

$ws.Range("C129").Value2 = 632760
$ws.Range("E129").Value2 = 3417328694

$ws.Range("C130").Value2 = 4227
$ws.Range("E130").Value2 = 139050789

$ws.Range("C132").Value2 = 585038
$ws.Range("E132").Value2 = 3441943057

$ws.Range("C136").Value2 = 26632
$ws.Range("E136").Value2 = 141844738

$ws.Range("C144").Value2 = 24513
$ws.Range("E144").Value2 = 88165149

$ws.Range("C151").Value2 = 39269
$ws.Range("E151").Value2 = 59810882

$ws.Range("C154").Value2 = 17975
$ws.Range("E154").Value2 = 69519959

$ws.Range("C157").Value2 = 630
$ws.Range("E157").Value2 = 1388087

$ws.Range("C171").Value2 = 95810
$ws.Range("E171").Value2 = 490330707

$ws.Range("C196").Value2 = 595464
$ws.Range("E196").Value2 = 983955743

$ws.Range("C215").Value2 = 230236
$ws.Range("E215").Value2 = 408662557

$ws.Range("C221").Value2 = 135467
$ws.Range("E221").Value2 = 681713588

$ws.Range("C229").Value2 = 612510
$ws.Range("E229").Value2 = 1040628728

$ws.Range("C237").Value2 = 283249
$ws.Range("E237").Value2 = 1437577049

$ws.Range("C240").Value2 = 205848
$ws.Range("E240").Value2 = 1066948938
